$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.115.90"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.309.29"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.86%  "
$ws.Range("D16").Value = "2.652.86"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "2.304.55"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "42.069.47"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  -6.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0892"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.04%  "
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.60%  "
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.231"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.48%  "
$ws.Range("E51").Value = "  -0.25%  "
